$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 132 (shifts existing rows 132-208 down to 133-209)
$ws.Rows(132).Insert()

# Populate the newly inserted row 132 with the new data record
$ws.Cells.Item(132, 1).Value = 7
$ws.Cells.Item(132, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(132, 3).Value = "Ñuble"
$ws.Cells.Item(132, 4).Value = 44606
$ws.Cells.Item(132, 5).Value = 16
$ws.Cells.Item(132, 6).Value = 100112009
$ws.Cells.Item(132, 7).Value = "Acelga"
$ws.Cells.Item(132, 8).Value = "Sin especificar"
$ws.Cells.Item(132, 9).Value = "Primera"
$ws.Cells.Item(132, 10).Value = 60
$ws.Cells.Item(132, 11).Value = 400
$ws.Cells.Item(132, 12).Value = 450
$ws.Cells.Item(132, 13).Value = 425
$ws.Cells.Item(132, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(132, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(132, 16).Value = 425
$ws.Cells.Item(132, 17).Value = 1
$ws.Cells.Item(132, 18).Value = "Hortaliza"
